$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.987.79'
$ws.Range('E2').Value = '  -0.68%  '
$ws.Range('D3').Value = '1.761.40'
$ws.Range('E3').Value = '  -3.39%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  -0.22%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '337.15'
$ws.Range('E5').Value = '  -1.02%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9988'
$ws.Range('E6').Value = '  -0.41%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3777'
$ws.Range('E7').Value = '  -4.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3356'
$ws.Range('E8').Value = '  -4.24%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '45.71'
$ws.Range('E9').Value = '  -5.03%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.122'
$ws.Range('E10').Value = '  -6.75%  '
$ws.Range('E11').Value = '  -5.35%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.001'
$ws.Range('E12').Value = '  -0.18%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '22.33'
$ws.Range('E13').Value = '  +0.23%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.195'
$ws.Range('E14').Value = '  -5.46%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.188'
$ws.Range('E15').Value = '  -0.44%  '
$ws.Range('D16').Value = '1.759.24'
$ws.Range('E16').Value = '  -3.68%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001053'
$ws.Range('E17').Value = '  -5.06%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06583'
$ws.Range('E18').Value = '  -2.27%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '80.36'
$ws.Range('E19').Value = '  -6.10%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.9994'
$ws.Range('E20').Value = '  -0.32%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '16.99'
$ws.Range('E21').Value = '  -5.05%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.286'
$ws.Range('E22').Value = '  -5.29%  '
$ws.Range('D23').Value = '27.998.12'
$ws.Range('E23').Value = '  -0.66%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.73'
$ws.Range('E24').Value = '  -7.47%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.370'
$ws.Range('E25').Value = '  -2.00%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '152.63'
$ws.Range('E26').Value = '  -1.75%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.86'
$ws.Range('E27').Value = '  -7.64%  '
$ws.Range('E28').Value = '  -8.74%  '
$ws.Range('D29').Value = '1.958.49'
$ws.Range('E29').Value = '  -3.70%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '131.88'
$ws.Range('E30').Value = '  -3.33%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.254'
$ws.Range('E31').Value = '  -16.26%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.015'
$ws.Range('E32').Value = '  -0.98%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.807'
$ws.Range('E33').Value = '  -6.44%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08793'
$ws.Range('E34').Value = '  -0.63%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '12.28'
$ws.Range('E35').Value = '  -6.78%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02342'
$ws.Range('E36').Value = '  -4.10%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.6604'
$ws.Range('E37').Value = '  -5.43%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.06175'
$ws.Range('E38').Value = '  -6.34%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.153'
$ws.Range('E39').Value = '  -6.75%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.2115'
$ws.Range('E40').Value = '  -5.41%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.212'
$ws.Range('E41').Value = '  -4.55%  '
$ws.Range('E42').Value = '  -10.54%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.018'
$ws.Range('E43').Value = '  -6.30%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9986'
$ws.Range('E44').Value = '  -0.40%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.74'
$ws.Range('E45').Value = '  -6.57%  '
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.6059'
$ws.Range('E46').Value = '  -6.95%  '
$ws.Range('B47').Value = 'PancakeSwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.822'
$ws.Range('E47').Value = '  -1.71%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '129.86'
$ws.Range('E48').Value = '  -1.57%  '
$ws.Range('E49').Value = '  -7.54%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.182'
$ws.Range('E50').Value = '  +1.63%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.07165'
$ws.Range('E51').Value = '  -0.91%  '
